## ---------------------------------------------------------------------
## Adds a new "robots" worksheet (a small metadata table about the
## individual robot / end-effector models used elsewhere in the
## workbook) right before the "analysis" sheet, and makes it the active
## sheet -- matching the author's commit "added meta information to
## random dataset".
## ---------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$analysisSheet = $wb.Worksheets.Item("analysis")

# New sheet gets inserted directly in front of "analysis" (i.e. right
# after "information"), matching sheets order: ... information, robots, analysis
$ws = $wb.Worksheets.Add($analysisSheet)
$ws.Name = "robots"

# --- column widths ------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 43.86
$ws.Columns.Item(2).ColumnWidth = 27.71
$ws.Columns.Item(3).ColumnWidth = 20.71

# --- header row ----------------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Manufacturer"
$ws.Range("A1:C1").Font.Bold = $true

# --- body rows -------------------------------------------------------
# Cells are written in an order that mirrors how the table was actually
# filled in (column A down for a couple of rows, then back to fill B/C),
# so newly-introduced shared strings land in the same order as the
# original edit.
$ws.Range("A2").Value = "UR5e"
$ws.Range("A3").Value = "Allegro hand - left"
$ws.Range("A4").Value = "Allegro hand - right"

$ws.Range("B2").Value = "robotic arm"
$ws.Range("B3").Value = "end effector"
$ws.Range("B4").Value = "end effector"

$ws.Range("C2").Value = "Universal Robots"
$ws.Range("C3").Value = "Wonik Robotics"
$ws.Range("C4").Value = "Wonik Robotics"

$ws.Range("A5").Value = "Atlas - convex hull"
$ws.Range("B5").Value = "humanoid robot"
$ws.Range("C5").Value = "Boston Dynamics"

$ws.Range("A6").Value = "Atlas - minimal contact"
$ws.Range("B6").Value = "humanoid robot"
$ws.Range("C6").Value = "Boston Dynamics"

$ws.Range("A7").Value = "3F gripper"
$ws.Range("B7").Value = "end effector"
$ws.Range("C7").Value = "Robotiq"

$ws.Range("A8").Value = "3F gripper - articulated"
$ws.Range("B8").Value = "end effector"
$ws.Range("C8").Value = "Robotiq"

$ws.Range("A9").Value = "3F gripper - tendons"
$ws.Range("B9").Value = "end effector"
$ws.Range("C9").Value = "Robotiq"

# trailing note, styled like an inline code / literal value
$ws.Range("A10").Value = "Kuk"
$ws.Range("A10").Font.Name = "Consolas"
$ws.Range("A10").Font.Family = 3
$ws.Range("A10").Font.Color = 7901646
$ws.Range("A10").VerticalAlignment = -4108

# --- view state: make "robots" the selected / active sheet ----------
$ws.Activate()
$excel.ActiveWindow.Zoom = 175
$ws.Range("A9").Select()

Write-Host "robots sheet added"
